# Miami.xlsx — "Updated the Excel files for the first 4 centers so that the
# markers are more straightforward to parse."
#
# The "Comp controls" sheet's Marker column (B4:B19) previously stored only
# the fluorophore/channel code (e.g. "FITC", "PE-A", "PE Cy7", ...). Several
# of those codes repeat (five rows all said "PE Cy7"), and a handful of
# marker names ("CD127 Alexa 647", "CD8 APC-H7", "CD45RO APC-H7",
# "CD3+19+20 APC-H7") were stuffed into what should have been a channel-only
# column. The fix: every Marker cell becomes "<name>:<channel>" so the
# antibody/marker and its detection channel are both visible and easy to
# parse.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Comp controls")

$ws.Range("B4").Value  = "LIVE GREEN:FITC"
$ws.Range("B5").Value  = "CD197:PE-A"
$ws.Range("B6").Value  = "CD4:PerCP-Cy5-5-A"
$ws.Range("B7").Value  = "CD45RA:PE Cy7"
$ws.Range("B8").Value  = "CD194:PE Cy7"
$ws.Range("B9").Value  = "CD27:PE Cy7"
$ws.Range("B10").Value = "CD11c:PE Cy7"
$ws.Range("B11").Value = "CD196:PE Cy7"
$ws.Range("B12").Value = "CD38:APC-A"
$ws.Range("B13").Value = "CD127:Alexa 647"
$ws.Range("B14").Value = "CD8:APC-H7"
$ws.Range("B15").Value = "CD45RO:APC-H7"
$ws.Range("B16").Value = "CD20:APC-Cy7-A"
$ws.Range("B17").Value = "CD3+19+20:APC-H7"
$ws.Range("B18").Value = "CD3:Pacific Blue-A"
$ws.Range("B19").Value = "HLA-DR:AmCyan-A"

# Bring the "Exp samples" sheet's zoom in line with "Comp controls" (150%),
# then leave the workbook focused back on "Comp controls" with the
# selection parked just below the filled-in table, matching the saved
# view state of the edited workbook.
$wsExp = $wb.Worksheets.Item("Exp samples")
$wsExp.Activate()
$excel.ActiveWindow.Zoom = 150

$ws.Activate()
[void]$ws.Range("B20").Select()
